# REVER_DailyTracker_MONISHA.xlsx — "Add files via upload" edit
#
# Target sheet: NOV-2020 (3rd tab, index 3 / 1-based) which is the
# tabSelected sheet in the workbook.
#
# Summary of the edit (rows are 1-based as in the sheet):
#   * Row 7  (D7)  : task text updated to mention nMVAR_Export_bat
#   * Row 8  (D8)  : "Week off" cell keeps its look (font/border) -- no
#                    value change, just kept as-is
#   * Row 10 (D10) : task text updated to mention nMVAR_Export_bat
#   * Row 10 (E10) : % of completion switched from a free-text note to a
#                    real 100% numeric value
#   * Row 11 (D11) : task text replaced with "nMVAR_Import_bat"
#   * Row 11 (F11) : status switched from COMPLETED to WIP
#   * Row 11 (G11) : comment cell border normalised to match the other
#                    WIP rows above it
#   * Row 11        : row height reset back to the sheet's default
#   * Row 12        : previously-blank placeholder row filled in with a
#                    new day's entry (nMVAR_Import_bat testing)
#   * Selection     : active cell moved to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Row 7 & 10: update the task descriptions (same style, new text)
# ---------------------------------------------------------------------
$ws.Range("D7").Value2  = "1) nMVAR_QA    2)nMVAR_Export_bat        "
$ws.Range("D10").Value2 = "1) nMVAR_QA                                                                     2)nMVAR_Export_bat"

# Row 10: % of completion is now a plain 100% number instead of free text
$ws.Range("E10").Value2 = 1

# ---------------------------------------------------------------------
# Row 11: new task text, status flips from COMPLETED back to WIP
# ---------------------------------------------------------------------
$ws.Range("D11").Value2 = "nMVAR_Import_bat"

# Re-style F11 like the other WIP rows (F3:F7) before overwriting its text
$ws.Range("F3").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F11").Value2 = "WIP"

# Re-style G11 (comment cell) like G10/G3, which use the plain 4-side
# border instead of the odd "no right edge" border it had before
$ws.Range("G10").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 11 no longer needs the tall custom row height
$ws.Rows.Item(11).EntireRow.AutoFit()

# ---------------------------------------------------------------------
# Row 12: fill in the previously-empty placeholder row with a new entry,
# reusing row 11's (now-fixed) formatting across the row
# ---------------------------------------------------------------------
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A12").Value2 = 11
$ws.Range("B12").Value2 = 44146
$ws.Range("C12").Value2 = "nMVAR "
$ws.Range("D12").Value2 = "nMVAR_Import_bat testing with import file"
$ws.Range("E12").Value2 = 0.7
$ws.Range("F12").Value2 = "WIP"

# ---------------------------------------------------------------------
# Selection moves to D11 (matches the saved cursor position in the file)
# ---------------------------------------------------------------------
$ws.Range("D11").Select()
